# Regenerate the localization-status handoff report:
#  - bump the two "latest" timestamps that changed because the report
#    was regenerated a few seconds later
#  - mark the Priority column ("ht" = hotfix/high-priority handoff) for
#    the rows that are now flagged, on both language sheets

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows whose handoff/handback bookkeeping changed in this run.
$rows = @(7, 8, 9, 10, 13, 14)

# Overview!G<row> = "Latest HO Xliff Generate Date" -> new timestamp.
# de-de!H<row> uses the very same timestamp string.
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-26 04:21:31"
    $wsDeDe.Range("H$r").Value     = "2016-08-26 04:21:31"
}

# zh-cn!H<row> = "Latest Handoff Datetime" -> its own new timestamp.
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-26 04:21:26"
}

# zh-cn!E<row> / de-de!E<row> = "Priority" -> now flagged "ht".
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
